# edit.ps1
# Applies the Case C.3 Ransomware Analysis commit:
#   1. Drop the stray _GoBack bookmark after
#      "It wasn’t stopped but the attack happened rarely "
#   2. Add the student's answer to the (previously empty) paragraph that
#      follows "How could you be affected by a ransomware attack?"
#   3. Add the final answer text (with a lastRenderedPageBreak run) to the
#      last (previously empty) paragraph, and move the _GoBack bookmark
#      there, followed by a trailing space run.
#
# Because this runtime coalesces adjacent same-formatted runs back into a
# single <w:r> when text is written via Range.Text, the runs here are
# built via Range.InsertXML with an explicit WordProcessingML package so
# the exact run/bookmark layout from the target document is reproduced.

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, [string]$pattern, [int]$fallbackIndex) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -match $pattern) {
            return $p
        }
    }
    return $doc.Paragraphs($fallbackIndex)
}

# ---------------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark from the "It wasn’t stopped..." answer
# ---------------------------------------------------------------------
$pStopped = Find-ParagraphByText $d "It wasn.?t stopped but the attack happened rarely" 40

$xmlStopped = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="0A1EF82C" w14:textId="124BCD06" w:rsidR="00E10B53" w:rsidRDefault="008C6CE9" w:rsidP="00E10B53">' +
  '<w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="360"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">It wasn’t stopped but the attack happened rarely </w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$pStopped.Range.InsertXML($xmlStopped)

# ---------------------------------------------------------------------
# 2) Fill in the answer right after
#    "How could you be affected by a ransomware attack?"
# ---------------------------------------------------------------------
$pRansomAnswer = Find-ParagraphByText $d "How could you be affected by a" 41
$pAnswerIndex = $pRansomAnswer.Range.Information(3)  # wdActiveEndAdjustedPageNumber unused; placeholder
$pEmpty1 = $d.Paragraphs($pRansomAnswer.Index + 1)

$xmlEmail = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="2386703E" w14:textId="77777777" w:rsidR="00E10B53" w:rsidRDefault="00E10B53" w:rsidP="00E10B53">' +
  '<w:pPr><w:pStyle w:val="NoSpacing"/><w:ind w:left="360"/></w:pPr>' +
  '<w:r><w:t>If I get an email that can hack m</w:t></w:r>' +
  '<w:r><w:t>y computer and ask for a ransom</w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$pEmpty1.Range.InsertXML($xmlEmail)

# ---------------------------------------------------------------------
# 3) Fill in the final (last) paragraph with the closing answer, the
#    relocated _GoBack bookmark, and a trailing space run.
# ---------------------------------------------------------------------
$pLastIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs($pLastIndex)
$lastStart = $pLast.Range.Start
$rLast = $d.Range($lastStart, $lastStart)

$xmlDoNot = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
  '<w:body>' +
  '<w:p w14:paraId="417C7537" w14:textId="77777777" w:rsidR="00E10B53" w:rsidRPr="00E10B53" w:rsidRDefault="00E10B53" w:rsidP="00E10B53">' +
  '<w:pPr><w:pStyle w:val="NoSpacing"/></w:pPr>' +
  '<w:r><w:lastRenderedPageBreak/><w:t>Do not open email or go to website which say free stuff or very good discount</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

[void]$rLast.InsertXML($xmlDoNot)

Write-Host "Applied ransomware-analysis edits. Final paragraph count:" $d.Paragraphs.Count
